$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("data")
$ws.Activate()

# Header for new column C
$ws.Range("C1").Value = "T"

# New column C data (temperature values) aligned with existing A/B rows
$ws.Range("C2").Value = 100
$ws.Range("C3").Value = 96.4
$ws.Range("C4").Value = 93.5
$ws.Range("C5").Value = 91.2
$ws.Range("C6").Value = 89.3
$ws.Range("C7").Value = 87.7
$ws.Range("C8").Value = 84.4
$ws.Range("C9").Value = 81.7
$ws.Range("C10").Value = 78
$ws.Range("C11").Value = 75.3
$ws.Range("C12").Value = 73.1
$ws.Range("C13").Value = 71.2
$ws.Range("C14").Value = 69.3
$ws.Range("C15").Value = 67.6
$ws.Range("C16").Value = 66
$ws.Range("C17").Value = 65
$ws.Range("C18").Value = 64.5

# Move the active selection to C19 (cell right below the new data), matching
# the saved workbook state in the target file.
$ws.Range("C19").Select()
